$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New result row (row 4) - attention result with 100 epochs
$ws.Range("A4").Value = "InceptV3_LSTM_Attention_CustEmbedding"
$ws.Range("B4").Value = "InceptionV3"
$ws.Range("C4").Value = "LSTM"
$ws.Range("E4").Value = 0.3015
$ws.Range("F4").Value = 0.3015
$ws.Range("G4").Value = 0.5192
$ws.Range("H4").Value = 0.5299

# G4/H4 get a wrapped-text style (new cellXfs entry), matching the rest of
# the "TRAIN" columns style but with word wrap enabled.
$ws.Range("G4:H4").HorizontalAlignment = 1
$ws.Range("G4:H4").WrapText = $true

# Column A got a bit wider to fit the new, longer model name (stored width
# 51.96 once Excel/LibreOffice pads it; COM's ColumnWidth is pre-padding).
$ws.Columns.Item(1).ColumnWidth = 51.126666666666665

# Selection moved, as recorded at save time.
$null = $ws.Range("A20").Select()

Write-Output "done"
